# Auto-generated edit script: update cryptos price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.395.91"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "1.994.73"
$ws.Range("E3").Value = "  -5.83%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'330.68"
$ws.Range("E5").Value = "  -4.55%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4949"
$ws.Range("E7").Value = "  -4.85%  "
$ws.Range("D8").Value = "'0.4205"
$ws.Range("E8").Value = "  -5.76%  "
$ws.Range("D9").Value = "'52.23"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "'0.08866"
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("D11").Value = "'1.118"
$ws.Range("E11").Value = "  -5.51%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'23.28"
$ws.Range("E12").Value = "  -8.39%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "2.011.52"
$ws.Range("E13").Value = "  -5.49%  "
$ws.Range("D14").Value = "'8.037"
$ws.Range("E14").Value = "  -6.68%  "
$ws.Range("D15").Value = "'6.510"
$ws.Range("E15").Value = "  -6.34%  "
$ws.Range("D16").Value = "'96.34"
$ws.Range("E16").Value = "  -6.09%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  -5.21%  "
$ws.Range("D19").Value = "'0.06618"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'19.72"
$ws.Range("E20").Value = "  -8.61%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'5.966"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").Value = "29.438.11"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").Value = "'11.85"
$ws.Range("E24").Value = "  -6.94%  "
$ws.Range("D25").Value = "'2.285"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'157.75"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").Value = "'6.599"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.55"
$ws.Range("E28").Value = "  -7.23%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.350"
$ws.Range("E29").Value = "  -7.37%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'127.51"
$ws.Range("E30").Value = "  -4.86%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.056"
$ws.Range("E31").Value = "  -8.35%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09941"
$ws.Range("E32").Value = "  -5.90%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.566"
$ws.Range("E33").Value = "  -12.57%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.860"
$ws.Range("E34").Value = "  -6.55%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.780"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.597"
$ws.Range("E36").Value = "  -10.95%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02454"
$ws.Range("E37").Value = "  -6.71%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06370"
$ws.Range("E38").Value = "  -7.28%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.285"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.78"
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6498"
$ws.Range("E41").Value = "  -8.35%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2071"
$ws.Range("E42").Value = "  -7.95%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.007"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6334"
$ws.Range("E44").Value = "  -7.63%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.47"
$ws.Range("E45").Value = "  -8.48%  "
$ws.Range("D46").Value = "'2.207"
$ws.Range("E46").Value = "  -7.59%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.267"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.542"
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000331"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'1.151"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06988"
$ws.Range("E51").Value = "  -2.59%  "
